$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers move from row 3 to row 1 (A1 already holds "ID Alternatif") ---
$ws.Range("B1").Value = "Nama Alternatif"
$ws.Range("C1").Value = "Nilai Raport (C1)"
$ws.Range("D1").Value = "Nilai Etika (C2)"
$ws.Range("E1").Value = "Nilai Kehadiran (C3)"
$ws.Range("F1").Value = "Nilai Ekstrakulikuler (C4)"

# --- Data row moves from row 4 to row 2 ---
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = "cek"

# --- Remove the old header row (row 3) and the old dropdown-source rows (4-20) ---
$ws.Range("A3:F20").ClearContents()

# --- Replace the old helper-column dropdown source ranges with real dropdown
#     list data validations on the row-2 cells ---
$ws.Range("C2").Validation.Add(3, 1, 1, """95-100,90-94,85-89,<=84""")
$ws.Range("C2").Validation.IgnoreBlank = $false
$ws.Range("C2").Validation.ShowInput = $false
$ws.Range("C2").Validation.ShowError = $false
$ws.Range("C2").Validation.ErrorTitle = "Input error"
$ws.Range("C2").Validation.ErrorMessage = "Please select a value from the list."
$ws.Range("C2").Validation.InputTitle = "Pick from list"
$ws.Range("C2").Validation.InputMessage = "Please pick a value from the dropdown list."

$ws.Range("D2").Validation.Add(3, 1, 1, """sangat baik,baik,cukup baik,kurang baik""")
$ws.Range("D2").Validation.IgnoreBlank = $false
$ws.Range("D2").Validation.ShowInput = $false
$ws.Range("D2").Validation.ShowError = $false
$ws.Range("D2").Validation.ErrorTitle = "Input error"
$ws.Range("D2").Validation.ErrorMessage = "Please select a value from the list."
$ws.Range("D2").Validation.InputTitle = "Pick from list"
$ws.Range("D2").Validation.InputMessage = "Please pick a value from the dropdown list."

$ws.Range("E2").Validation.Add(3, 1, 1, """selalu hadir,cukup hadir,jarang hadir,izin""")
$ws.Range("E2").Validation.IgnoreBlank = $false
$ws.Range("E2").Validation.ShowInput = $false
$ws.Range("E2").Validation.ShowError = $false
$ws.Range("E2").Validation.ErrorTitle = "Input error"
$ws.Range("E2").Validation.ErrorMessage = "Please select a value from the list."
$ws.Range("E2").Validation.InputTitle = "Pick from list"
$ws.Range("E2").Validation.InputMessage = "Please pick a value from the dropdown list."

$ws.Range("F2").Validation.Add(3, 1, 1, """sangat aktif,aktif,cukup aktif,kurang aktif""")
$ws.Range("F2").Validation.IgnoreBlank = $false
$ws.Range("F2").Validation.ShowInput = $false
$ws.Range("F2").Validation.ShowError = $false
$ws.Range("F2").Validation.ErrorTitle = "Input error"
$ws.Range("F2").Validation.ErrorMessage = "Please select a value from the list."
$ws.Range("F2").Validation.InputTitle = "Pick from list"
$ws.Range("F2").Validation.InputMessage = "Please pick a value from the dropdown list."
